$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.993.57'

$ws.Range("D3").Value = '3.745.78'
$ws.Range("E3").Value = '  -0.49%  '

$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").Value = "'602.13"
$ws.Range("E5").Value = '  -0.18%  '

$ws.Range("D6").Value = "'165.55"
$ws.Range("E6").Value = '  -2.40%  '

$ws.Range("D7").Value = '3.744.41'
$ws.Range("E7").Value = '  -0.44%  '

$ws.Range("E8").Value = '  +0.24%  '

$ws.Range("D9").Value = "'0.538"
$ws.Range("E9").Value = '  +0.24%  '

$ws.Range("D10").Value = "'0.171"
$ws.Range("E10").Value = '  +4.13%  '

$ws.Range("D11").Value = "'6.36"
$ws.Range("E11").Value = '  +0.27%  '

$ws.Range("E12").Value = '  -1.33%  '

$ws.Range("D13").Value = "'37.61"
$ws.Range("E13").Value = '  -2.30%  '

$ws.Range("E14").Value = '  +0.25%  '

$ws.Range("D15").Value = '4.372.56'
$ws.Range("E15").Value = '  -0.49%  '

$ws.Range("D16").Value = '3.758.08'
$ws.Range("E16").Value = '  +0.15%  '

$ws.Range("D17").Value = '68.924.91'
$ws.Range("E17").Value = '  +0.17%  '

$ws.Range("E18").Value = '  +1.27%  '

$ws.Range("D19").Value = "'17.80"
$ws.Range("E19").Value = '  +3.63%  '

$ws.Range("E20").Value = '  -0.98%  '

$ws.Range("D21").Value = "'11.29"
$ws.Range("E21").Value = '  +4.35%  '

$ws.Range("D22").Value = "'489.46"
$ws.Range("E22").Value = '  -1.71%  '

$ws.Range("D23").Value = "'0.723"
$ws.Range("E23").Value = '  -1.13%  '

$ws.Range("E24").Value = '  +2.36%  '

$ws.Range("D25").Value = "'84.51"
$ws.Range("E25").Value = '  -1.57%  '

$ws.Range("E26").Value = '  -2.96%  '

$ws.Range("D27").Value = "'12.26"
$ws.Range("E27").Value = '  -1.06%  '

$ws.Range("E28").Value = '  -2.17%  '

$ws.Range("E29").Value = '  +0.08%  '

$ws.Range("D30").Value = "'2.96"
$ws.Range("E30").Value = '  -0.72%  '

$ws.Range("E31").Value = '  +1.89%  '

$ws.Range("E32").Value = '  -5.01%  '

$ws.Range("D33").Value = '3.890.96'
$ws.Range("E33").Value = '  -0.54%  '

$ws.Range("D34").Value = "'31.54"
$ws.Range("E34").Value = '  -2.16%  '

$ws.Range("D35").Value = '3.684.97'

$ws.Range("D36").Value = "'0.107"
$ws.Range("E36").Value = '  -1.09%  '

$ws.Range("D37").Value = "'5.92"
$ws.Range("E37").Value = '  +0.97%  '

$ws.Range("E38").Value = '  -0.98%  '

$ws.Range("E39").Value = '  +4.25%  '

$ws.Range("D40").Value = "'0.999"
$ws.Range("E40").Value = '  -0.03%  '

$ws.Range("D41").Value = "'0.324"
$ws.Range("E41").Value = '  -0.74%  '

$ws.Range("D42").Value = "'3.07"
$ws.Range("E42").Value = '  +6.90%  '

$ws.Range("D43").Value = "'48.52"
$ws.Range("E43").Value = '  -0.61%  '

$ws.Range("E44").Value = '  +0.46%  '

$ws.Range("D45").Value = "'423.56"
$ws.Range("E45").Value = '  -4.80%  '

$ws.Range("B46").Value = 'USDe'
$ws.Range("C46").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D46").Value = "'1.00"
$ws.Range("E46").Value = '  +0.03%  '

$ws.Range("B47").Value = 'Cosmos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D47").Value = "'8.41"
$ws.Range("E47").Value = '  -1.21%  '

$ws.Range("B48").Value = 'Arweave'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D48").Value = "'40.00"
$ws.Range("E48").Value = '  -1.82%  '

$ws.Range("B49").Value = 'Monero'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D49").Value = "'141.78"
$ws.Range("E49").Value = '  -0.11%  '

$ws.Range("B50").Value = 'ONDO'
$ws.Range("C50").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D50").Value = "'1.30"
$ws.Range("E50").Value = '  +7.26%  '

$ws.Range("D51").Value = '2.780.28'
$ws.Range("E51").Value = '  -2.21%  '
